# Update the "想去人数" (wishlist count) and "最低票价" (min price) figures
# for the events shown on both the "展览" and "全部类型" worksheets, which
# carry duplicate copies of the same rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1465
    $ws.Range("G2").Value = 39.9

    $ws.Range("F3").Value = 3077

    $ws.Range("F4").Value = 41

    $ws.Range("F5").Value = 646
}

$wb.Save()
